$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M_MUT_and_WT_M_P60_CORT")

# --- Formulas for columns B, C, D (DESeq2/Limma/EdgeR totals) ---
$ws.Range("B2").Formula = "=E2+H2+I2+K2"
$ws.Range("C2").Formula = "=F2+H2+J2+K2"
$ws.Range("D2").Formula = "=G2+J2+I2+K2"

$ws.Range("B3:B15").Formula = "=E3+H3+I3+K3"
$ws.Range("C3:C15").Formula = "=F3+H3+J3+K3"
$ws.Range("D3:D15").Formula = "=G3+J3+I3+K3"

# --- Column K formatting (new cell styles) ---
# K1 header cell: bold + (re-stamped) fill + center/center
$ws.Range("K1").Interior.Color = 16777215

# K2:K7: font explicitly applied + fill + center/center
$ws.Range("K2:K7").Interior.Color = 16777215
$ws.Range("K2:K7").Font.ThemeColor = 1

# K8:K15: fill + center/center (no explicit font)
$ws.Range("K8:K10").Interior.Color = 16777215
$ws.Range("K11:K12").Interior.Color = 16777215
$ws.Range("K11:K12").HorizontalAlignment = -4108
$ws.Range("K11:K12").VerticalAlignment = -4108
$ws.Range("K13:K15").Interior.Color = 16777215

# --- New row 16, column K only, fill-only style, no alignment ---
$ws.Range("K16").Interior.Color = 16777215

# --- Column K width (bestFit-like) ---
$ws.Columns.Item(11).ColumnWidth = 10.25

# --- Selection / active cell ---
$ws.Range("I8").Select()
